# Fill in the "Row" column (C) on Sheet1 with the rank/order values
# corresponding to each weekday's earnings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2" = 2
    "C3" = 4
    "C4" = 1
    "C5" = 5
    "C6" = 3
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
